# Apply the edit described by the commit:
# "Mise en place d'un OnItemClick listener sur la ListView permettant d'appuyer sur
#  les items et d'un OnTouch listener qui desactive le scroll du parent."
#
# On the "Iteration #1" worksheet, two new iteration-log entries are added
# (rows 18 and 19), one blank placeholder row is removed from the bottom of the
# data-entry block, and the TOTAL formula / merged cells shift accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #1")

# --- Row 18: turn the first still-empty row into a real log entry ---
# Copy the date formatting from A17 (keeps numFmtId=14 / style 12 intact)
$ws.Cells.Item(17, 1).Copy($ws.Cells.Item(18, 1))
$ws.Cells.Item(18, 1).Value = 42766

$ws.Cells.Item(18, 2).Value = "Le sous-menu est désormais une listview avec un adapter."
$ws.Cells.Item(18, 3).Value = 2

# --- Row 19: a brand-new log entry with a longer, wrapped description ---
$ws.Cells.Item(17, 1).Copy($ws.Cells.Item(19, 1))
$ws.Cells.Item(19, 1).Value = 42771
$ws.Cells.Item(19, 1).VerticalAlignment = -4108  # xlVAlignCenter

$ws.Cells.Item(19, 2).Value = "Click listener et touch listener sur la listview pour activer le scroll et appuyer sur les items."
$ws.Cells.Item(19, 2).WrapText = $true

$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 3).HorizontalAlignment = -4152  # xlRight
$ws.Cells.Item(19, 3).VerticalAlignment = -4108    # xlCenter

# --- Remove one now-superfluous blank placeholder row from the data block ---
$ws.Rows.Item(36).Delete()

# --- Left-align the (still empty) date column for the remaining blank rows ---
$ws.Range("A20:A35").HorizontalAlignment = -4131  # xlLeft

Write-Host "Edit applied"
